# PWData.xlsx regression-suite refresh:
#  - swap the sample lease tenant (Virat Kohli -> Abdul Kalam) on the
#    "Leases" sheet, including the derived email address
#  - nudge the current selection on "Leases" from C4 to B5
#  - tidy column widths on both sheets to match the latest authoring pass

$wb = $excel.ActiveWorkbook

$wsProps  = $wb.Worksheets.Item("Properties")
$wsLeases = $wb.Worksheets.Item("Leases")

# --- Leases!A2:C2 - replace the sample tenant row -----------------------
$wsLeases.Range("A2").Value = "Abdul"
$wsLeases.Range("B2").Value = "Kalam"
$wsLeases.Range("C2").Value = "abdul.kalam@gmail.com"

# --- selection bookkeeping ------------------------------------------------
$wsLeases.Activate() | Out-Null
$wsLeases.Range("B5").Select() | Out-Null

# --- column width touch-up (values match the latest saved widths) --------
$wsProps.Columns.Item(1).ColumnWidth = 12.65
$wsProps.Columns.Item(2).ColumnWidth = 16.65

$wsLeases.Columns.Item(1).ColumnWidth = 13.3
$wsLeases.Columns.Item(3).ColumnWidth = 22.78
